# Updated symbol list on Thu Dec 22 23:50:32 UTC 2022 with GitHub Actions
# Refresh the Price column (column D) in the cryptos sheet with the latest
# quotes. Values are entered with a leading apostrophe so Excel keeps them
# as literal text (matching the sheet's existing text-formatted Price
# column) instead of silently reinterpreting/truncating them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'245.80"
$ws.Cells.Item(3, 4).Value = "'21.84"
$ws.Cells.Item(4, 4).Value = "'5.443"
$ws.Cells.Item(5, 4).Value = "'0.05770"
$ws.Cells.Item(6, 4).Value = "'3.414"
$ws.Cells.Item(7, 4).Value = "'6.326"
$ws.Cells.Item(8, 4).Value = "'0.8184"
$ws.Cells.Item(9, 4).Value = "'1.021"
$ws.Cells.Item(11, 4).Value = "'0.07299"
$ws.Cells.Item(12, 4).Value = "'0.03109"
$ws.Cells.Item(13, 4).Value = "'0.03091"
$ws.Cells.Item(14, 4).Value = "'4.140"
$ws.Cells.Item(15, 4).Value = "'0.09404"
$ws.Cells.Item(16, 4).Value = "'0.001598"
$ws.Cells.Item(17, 4).Value = "'0.04817"
$ws.Cells.Item(18, 4).Value = "'0.0005848"
$ws.Cells.Item(19, 4).Value = "'0.006249"
$ws.Cells.Item(20, 4).Value = "'0.004118"
$ws.Cells.Item(21, 4).Value = "'0.0009943"
$ws.Cells.Item(22, 4).Value = "'0.0001500"
$ws.Cells.Item(23, 4).Value = "'3.743"
$ws.Cells.Item(27, 4).Value = "'0.0003998"
$ws.Cells.Item(40, 4).Value = "'0.03887"
$ws.Cells.Item(41, 4).Value = "'0.006691"
$ws.Cells.Item(42, 4).Value = "'0.1070"
$ws.Cells.Item(44, 4).Value = "'0.006680"
$ws.Cells.Item(45, 4).Value = "'0.00005609"
$ws.Cells.Item(47, 4).Value = "'0.3899"
$ws.Cells.Item(50, 4).Value = "'0.01010"
